# "Averaged Data Table created"
#
# The sheet's Table1 originally had several columns (TPS, Energy Use per
# Transaction, Nakamoto Coefficient, % of nodes required to take over
# network, Strengths, Weaknesses) that were only partly populated. This
# fills the remaining blank data cells with "N/A" placeholders and makes
# sure each column carries a sensible number format:
#   - TPS / Energy Use per Transaction / Nakamoto Coefficient -> 0.00
#   - % of nodes required to take over network                -> 0%
#   - Strengths / Weaknesses                                   -> Text (@)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TPS (B), Energy Use per Transaction (C), Nakamoto Coefficient (D)
$ws.Range("B2:D10").NumberFormat = "0.00"
$ws.Range("B2:D10").Value = "N/A"

# % of nodes required to take over network (E) - rows 2:3 already had data,
# rows 4:10 were blank
$ws.Range("E2:E10").NumberFormat = "0%"
$ws.Range("E4:E10").Value = "N/A"

# Strengths (F), Weaknesses (G)
$ws.Range("F2:G10").NumberFormat = "@"
$ws.Range("F2:G10").Value = "N/A"

# Cursor position at save time
$null = $ws.Range("J16:K16").Select()
